$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric must be kept as text (matches source inlineStr type)
$textFormatCells = @("D6", "D7", "D11", "D16", "D17", "D19", "D22", "D23", "D24", "D26", "D27", "D28", "D29", "D30", "D31", "D33", "D34", "D35", "D36", "D38", "D40", "D41", "D42", "D45", "D46", "D47", "D48", "D49")
foreach ($cell in $textFormatCells) {
    $ws.Range($cell).NumberFormat = "@"
}

# Apply updated cell values
$ws.Range("D2").Value = "29.364.44"
$ws.Range("E2").Value = "  +0.24%  "
$ws.Range("D3").Value = "1.842.54"
$ws.Range("E3").Value = "  -0.08%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("E5").Value = "  -0.27%  "
$ws.Range("D6").Value = "0.6308"
$ws.Range("E6").Value = "  +0.42%  "
$ws.Range("D7").Value = "0.9996"
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("E8").Value = "  +0.27%  "
$ws.Range("E9").Value = "  +0.17%  "
$ws.Range("E10").Value = "  +2.84%  "
$ws.Range("D11").Value = "0.07726"
$ws.Range("D12").Value = "1.843.81"
$ws.Range("E12").Value = "  +0.01%  "
$ws.Range("E13").Value = "  -0.15%  "
$ws.Range("E14").Value = "  -0.03%  "
$ws.Range("E15").Value = "  +2.20%  "
$ws.Range("D16").Value = "81.92"
$ws.Range("E16").Value = "  -0.18%  "
$ws.Range("D17").Value = "6.226"
$ws.Range("E17").Value = "  +1.43%  "
$ws.Range("D18").Value = "29.348.19"
$ws.Range("E18").Value = "  +0.23%  "
$ws.Range("D19").Value = "229.14"
$ws.Range("E19").Value = "  +0.53%  "
$ws.Range("E20").Value = "  +0.47%  "
$ws.Range("E21").Value = "  -0.13%  "
$ws.Range("D22").Value = "7.396"
$ws.Range("E22").Value = "  +0.14%  "
$ws.Range("D23").Value = "0.9999"
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("D24").Value = "158.23"
$ws.Range("E24").Value = "  -0.34%  "
$ws.Range("E25").Value = "  +1.64%  "
$ws.Range("D26").Value = "0.1355"
$ws.Range("E26").Value = "  -1.39%  "
$ws.Range("D27").Value = "17.47"
$ws.Range("E27").Value = "  -0.50%  "
$ws.Range("D28").Value = "0.06897"
$ws.Range("E28").Value = "  +10.10%  "
$ws.Range("D29").Value = "1.460"
$ws.Range("E29").Value = "  +4.59%  "
$ws.Range("D30").Value = "1.485"
$ws.Range("E30").Value = "  +0.84%  "
$ws.Range("D31").Value = "4.076"
$ws.Range("E31").Value = "  +0.72%  "
$ws.Range("E32").Value = "  -0.41%  "
$ws.Range("D33").Value = "1.833"
$ws.Range("E33").Value = "  +0.74%  "
$ws.Range("D34").Value = "1.141"
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("D35").Value = "0.7007"
$ws.Range("E35").Value = "  +1.26%  "
$ws.Range("D36").Value = "2.584"
$ws.Range("E36").Value = "  +0.06%  "
$ws.Range("E37").Value = "  +1.81%  "
$ws.Range("D38").Value = "2.817"
$ws.Range("E38").Value = "  -0.84%  "
$ws.Range("D39").Value = "1.238.03"
$ws.Range("E39").Value = "  -0.57%  "
$ws.Range("D40").Value = "6.816"
$ws.Range("E40").Value = "  +4.58%  "
$ws.Range("D41").Value = "0.9425"
$ws.Range("E41").Value = "  +3.71%  "
$ws.Range("D42").Value = "0.9987"
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("D43").Value = "1.993.37"
$ws.Range("E43").Value = "  -0.42%  "
$ws.Range("E44").Value = "  -0.41%  "
$ws.Range("D45").Value = "65.53"
$ws.Range("E45").Value = "  -1.02%  "
$ws.Range("D46").Value = "0.00000000119"
$ws.Range("E46").Value = "  +4.85%  "
$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").Value = "7.048"
$ws.Range("E47").Value = "  -0.13%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "1.716"
$ws.Range("E48").Value = "  +3.44%  "
$ws.Range("D49").Value = "9.003"
$ws.Range("E49").Value = "  -0.11%  "
$ws.Range("E50").Value = "  -1.48%  "
$ws.Range("E51").Value = "  -0.64%  "
